# Apply updated Betfair back/lay odds values to Sheet1 (data rows 2-25).
# Each statement mirrors one cell change from the source diff:
# $ws.Cells.Item(row, col).Value = newValue   (col numbers: F=6 ... AO=41)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 1.47  # F2: 1.46 -> 1.47
$ws.Cells.Item(2, 9).Value = 9.199999999999999  # I2: 8.199999999999999 -> 9.199999999999999
$ws.Cells.Item(2, 10).Value = 4.5  # J2: 4.1 -> 4.5
$ws.Cells.Item(2, 11).Value = 5.3  # K2: 6.2 -> 5.3
# Row 3
$ws.Cells.Item(3, 18).Value = 1.49  # R3: 1.48 -> 1.49
$ws.Cells.Item(3, 21).Value = 2.38  # U3: 2.36 -> 2.38
# Row 4
$ws.Cells.Item(4, 14).Value = 3.95  # N4: 3.9 -> 3.95
$ws.Cells.Item(4, 18).Value = 1.39  # R4: 1.37 -> 1.39
$ws.Cells.Item(4, 20).Value = 1.9  # T4: 1.91 -> 1.9
$ws.Cells.Item(4, 21).Value = 1.93  # U4: 1.92 -> 1.93
$ws.Cells.Item(4, 24).Value = 17.5  # X4: 19.5 -> 17.5
$ws.Cells.Item(4, 25).Value = 23  # Y4: 26 -> 23
$ws.Cells.Item(4, 28).Value = 10  # AB4: 10.5 -> 10
$ws.Cells.Item(4, 29).Value = 10.5  # AC4: 12 -> 10.5
$ws.Cells.Item(4, 37).Value = 20  # AK4: 21 -> 20
# Row 5
$ws.Cells.Item(5, 7).Value = 4.8  # G5: 4.9 -> 4.8
$ws.Cells.Item(5, 9).Value = 1.94  # I5: 1.95 -> 1.94
# Row 6
$ws.Cells.Item(6, 8).Value = 7.8  # H6: 8 -> 7.8
$ws.Cells.Item(6, 12).Value = 1.36  # L6: 1.43 -> 1.36
$ws.Cells.Item(6, 19).Value = 3.5  # S6: 3.85 -> 3.5
$ws.Cells.Item(6, 24).Value = 15.5  # X6: 1000 -> 15.5
$ws.Cells.Item(6, 28).Value = 7.6  # AB6: 6.8 -> 7.6
$ws.Cells.Item(6, 29).Value = 12.5  # AC6: 1000 -> 12.5
$ws.Cells.Item(6, 37).Value = 980  # AK6: 23 -> 980
$ws.Cells.Item(6, 40).Value = 12  # AN6: 1000 -> 12
# Row 7
$ws.Cells.Item(7, 30).Value = 970  # AD7: 980 -> 970
# Row 8
$ws.Cells.Item(8, 7).Value = 1.44  # G8: 1.45 -> 1.44
$ws.Cells.Item(8, 23).Value = 3.25  # W8: 3.2 -> 3.25
$ws.Cells.Item(8, 28).Value = 15  # AB8: 15.5 -> 15
$ws.Cells.Item(8, 29).Value = 16  # AC8: 17 -> 16
$ws.Cells.Item(8, 33).Value = 12.5  # AG8: 13 -> 12.5
$ws.Cells.Item(8, 36).Value = 16  # AJ8: 16.5 -> 16
$ws.Cells.Item(8, 37).Value = 16.5  # AK8: 17 -> 16.5
# Row 9
$ws.Cells.Item(9, 7).Value = 2.52  # G9: 2.66 -> 2.52
$ws.Cells.Item(9, 9).Value = 3.55  # I9: 3.95 -> 3.55
$ws.Cells.Item(9, 10).Value = 3.1  # J9: 3.4 -> 3.1
$ws.Cells.Item(9, 12).Value = 1.31  # L9: 1.01 -> 1.31
$ws.Cells.Item(9, 13).Value = 1.05  # M9: 1.01 -> 1.05
$ws.Cells.Item(9, 14).Value = 3.8  # N9: 1.98 -> 3.8
$ws.Cells.Item(9, 18).Value = 1.38  # R9: 1.32 -> 1.38
$ws.Cells.Item(9, 19).Value = 3  # S9: 2.66 -> 3
$ws.Cells.Item(9, 20).Value = 1.68  # T9: 1.01 -> 1.68
$ws.Cells.Item(9, 21).Value = 2.18  # U9: 1.01 -> 2.18
$ws.Cells.Item(9, 22).Value = 1.4  # V9: 1.35 -> 1.4
$ws.Cells.Item(9, 23).Value = 1.66  # W9: 1.62 -> 1.66
$ws.Cells.Item(9, 24).Value = 19.5  # X9: 1000 -> 19.5
$ws.Cells.Item(9, 25).Value = 17  # Y9: 1000 -> 17
$ws.Cells.Item(9, 26).Value = 29  # Z9: 1000 -> 29
$ws.Cells.Item(9, 27).Value = 70  # AA9: 1000 -> 70
$ws.Cells.Item(9, 28).Value = 13.5  # AB9: 1000 -> 13.5
$ws.Cells.Item(9, 29).Value = 8.800000000000001  # AC9: 1000 -> 8.800000000000001
$ws.Cells.Item(9, 30).Value = 17  # AD9: 1000 -> 17
$ws.Cells.Item(9, 31).Value = 44  # AE9: 1000 -> 44
$ws.Cells.Item(9, 32).Value = 19.5  # AF9: 1000 -> 19.5
$ws.Cells.Item(9, 33).Value = 14  # AG9: 1000 -> 14
$ws.Cells.Item(9, 34).Value = 20  # AH9: 1000 -> 20
$ws.Cells.Item(9, 35).Value = 55  # AI9: 1000 -> 55
$ws.Cells.Item(9, 36).Value = 38  # AJ9: 1000 -> 38
$ws.Cells.Item(9, 37).Value = 30  # AK9: 1000 -> 30
$ws.Cells.Item(9, 38).Value = 44  # AL9: 1000 -> 44
$ws.Cells.Item(9, 39).Value = 100  # AM9: 1000 -> 100
$ws.Cells.Item(9, 40).Value = 22  # AN9: 1000 -> 22
$ws.Cells.Item(9, 41).Value = 38  # AO9: 1000 -> 38
# Row 10
$ws.Cells.Item(10, 25).Value = 9.800000000000001  # Y10: 10 -> 9.800000000000001
# Row 11
$ws.Cells.Item(11, 7).Value = 6  # G11: 6.6 -> 6
$ws.Cells.Item(11, 8).Value = 1.78  # H11: 1.76 -> 1.78
$ws.Cells.Item(11, 9).Value = 1.95  # I11: 1.99 -> 1.95
$ws.Cells.Item(11, 10).Value = 3.4  # J11: 3.35 -> 3.4
$ws.Cells.Item(11, 11).Value = 3.9  # K11: 3.95 -> 3.9
$ws.Cells.Item(11, 12).Value = 1.42  # L11: 1.41 -> 1.42
$ws.Cells.Item(11, 13).Value = 1.09  # M11: 1.01 -> 1.09
$ws.Cells.Item(11, 18).Value = 1.25  # R11: 1.22 -> 1.25
$ws.Cells.Item(11, 19).Value = 4.1  # S11: 3.55 -> 4.1
$ws.Cells.Item(11, 20).Value = 1.98  # T11: 1.01 -> 1.98
$ws.Cells.Item(11, 21).Value = 1.8  # U11: 1.01 -> 1.8
$ws.Cells.Item(11, 22).Value = 2.04  # V11: 2 -> 2.04
$ws.Cells.Item(11, 28).Value = 1000  # AB11: 16 -> 1000
$ws.Cells.Item(11, 41).Value = 1000  # AO11: 17 -> 1000
# Row 12
$ws.Cells.Item(12, 7).Value = 2.3  # G12: 2.24 -> 2.3
$ws.Cells.Item(12, 9).Value = 4.8  # I12: 4.9 -> 4.8
$ws.Cells.Item(12, 12).Value = 1.45  # L12: 1.41 -> 1.45
$ws.Cells.Item(12, 19).Value = 3.75  # S12: 3.7 -> 3.75
$ws.Cells.Item(12, 23).Value = 1.78  # W12: 1.8 -> 1.78
$ws.Cells.Item(12, 26).Value = 980  # Z12: 32 -> 980
$ws.Cells.Item(12, 36).Value = 980  # AJ12: 32 -> 980
$ws.Cells.Item(12, 37).Value = 980  # AK12: 29 -> 980
$ws.Cells.Item(12, 40).Value = 980  # AN12: 24 -> 980
# Row 13
$ws.Cells.Item(13, 7).Value = 2.68  # G13: 2.62 -> 2.68
$ws.Cells.Item(13, 11).Value = 3.85  # K13: 3.55 -> 3.85
$ws.Cells.Item(13, 12).Value = 1.47  # L13: 1.43 -> 1.47
$ws.Cells.Item(13, 13).Value = 1.09  # M13: 1.08 -> 1.09
$ws.Cells.Item(13, 16).Value = 1.67  # P13: 1.66 -> 1.67
$ws.Cells.Item(13, 20).Value = 1.89  # T13: 1.79 -> 1.89
$ws.Cells.Item(13, 21).Value = 1.92  # U13: 1.82 -> 1.92
$ws.Cells.Item(13, 23).Value = 1.59  # W13: 1.61 -> 1.59
$ws.Cells.Item(13, 24).Value = 12  # X13: 12.5 -> 12
$ws.Cells.Item(13, 31).Value = 980  # AE13: 55 -> 980
# Row 14
$ws.Cells.Item(14, 6).Value = 3.45  # F14: 3.75 -> 3.45
$ws.Cells.Item(14, 7).Value = 5.3  # G14: 5.8 -> 5.3
$ws.Cells.Item(14, 8).Value = 1.81  # H14: 1.75 -> 1.81
$ws.Cells.Item(14, 9).Value = 2.02  # I14: 2.16 -> 2.02
$ws.Cells.Item(14, 10).Value = 3.45  # J14: 3.4 -> 3.45
$ws.Cells.Item(14, 11).Value = 5  # K14: 6.6 -> 5
$ws.Cells.Item(14, 14).Value = 2.14  # N14: 1.98 -> 2.14
$ws.Cells.Item(14, 16).Value = 2.14  # P14: 1.98 -> 2.14
$ws.Cells.Item(14, 17).Value = 1.61  # Q14: 1.6 -> 1.61
$ws.Cells.Item(14, 22).Value = 1.98  # V14: 1.86 -> 1.98
$ws.Cells.Item(14, 23).Value = 1.23  # W14: 1.21 -> 1.23
# Row 15
$ws.Cells.Item(15, 20).Value = 1.96  # T15: 1.81 -> 1.96
# Row 16
$ws.Cells.Item(16, 9).Value = 3.85  # I16: 3.9 -> 3.85
$ws.Cells.Item(16, 23).Value = 1.92  # W16: 1.93 -> 1.92
$ws.Cells.Item(16, 30).Value = 20  # AD16: 21 -> 20
# Row 17
$ws.Cells.Item(17, 7).Value = 4.5  # G17: 4.4 -> 4.5
$ws.Cells.Item(17, 12).Value = 1.28  # L17: 1.25 -> 1.28
$ws.Cells.Item(17, 17).Value = 1.57  # Q17: 1.56 -> 1.57
$ws.Cells.Item(17, 18).Value = 1.57  # R17: 1.58 -> 1.57
$ws.Cells.Item(17, 19).Value = 2.38  # S17: 2.4 -> 2.38
# Row 18
$ws.Cells.Item(18, 6).Value = 1.82  # F18: 1.71 -> 1.82
$ws.Cells.Item(18, 7).Value = 2.04  # G18: 2.16 -> 2.04
$ws.Cells.Item(18, 8).Value = 4.6  # H18: 4.3 -> 4.6
$ws.Cells.Item(18, 9).Value = 5.7  # I18: 6.4 -> 5.7
$ws.Cells.Item(18, 11).Value = 4.2  # K18: 5.4 -> 4.2
$ws.Cells.Item(18, 15).Value = 1.4  # O18: 1.39 -> 1.4
$ws.Cells.Item(18, 20).Value = 1.96  # T18: 1.01 -> 1.96
$ws.Cells.Item(18, 21).Value = 1.83  # U18: 1.01 -> 1.83
$ws.Cells.Item(18, 22).Value = 1.21  # V18: 1.19 -> 1.21
$ws.Cells.Item(18, 23).Value = 1.96  # W18: 1.86 -> 1.96
$ws.Cells.Item(18, 24).Value = 13.5  # X18: 1000 -> 13.5
$ws.Cells.Item(18, 26).Value = 44  # Z18: 1000 -> 44
$ws.Cells.Item(18, 27).Value = 160  # AA18: 1000 -> 160
$ws.Cells.Item(18, 28).Value = 8.6  # AB18: 1000 -> 8.6
$ws.Cells.Item(18, 29).Value = 9.199999999999999  # AC18: 1000 -> 9.199999999999999
$ws.Cells.Item(18, 30).Value = 24  # AD18: 1000 -> 24
$ws.Cells.Item(18, 31).Value = 90  # AE18: 1000 -> 90
$ws.Cells.Item(18, 32).Value = 12.5  # AF18: 1000 -> 12.5
$ws.Cells.Item(18, 33).Value = 12  # AG18: 1000 -> 12
$ws.Cells.Item(18, 34).Value = 26  # AH18: 1000 -> 26
$ws.Cells.Item(18, 35).Value = 100  # AI18: 1000 -> 100
$ws.Cells.Item(18, 36).Value = 26  # AJ18: 1000 -> 26
$ws.Cells.Item(18, 37).Value = 27  # AK18: 1000 -> 27
$ws.Cells.Item(18, 38).Value = 55  # AL18: 1000 -> 55
$ws.Cells.Item(18, 39).Value = 170  # AM18: 1000 -> 170
$ws.Cells.Item(18, 41).Value = 130  # AO18: 1000 -> 130
# Row 19
$ws.Cells.Item(19, 20).Value = 1.91  # T19: 2.06 -> 1.91
$ws.Cells.Item(19, 21).Value = 1.82  # U19: 1.71 -> 1.82
# Row 20
$ws.Cells.Item(20, 12).Value = 1.46  # L20: 1.01 -> 1.46
$ws.Cells.Item(20, 14).Value = 3.05  # N20: 2.76 -> 3.05
$ws.Cells.Item(20, 15).Value = 1.39  # O20: 1.37 -> 1.39
$ws.Cells.Item(20, 17).Value = 2.16  # Q20: 2.14 -> 2.16
$ws.Cells.Item(20, 19).Value = 3.55  # S20: 3.9 -> 3.55
$ws.Cells.Item(20, 20).Value = 1.84  # T20: 1.83 -> 1.84
$ws.Cells.Item(20, 21).Value = 1.98  # U20: 1.83 -> 1.98
# Row 22
$ws.Cells.Item(22, 21).Value = 2.78  # U22: 2.8 -> 2.78
$ws.Cells.Item(22, 27).Value = 75  # AA22: 80 -> 75
# Row 23
$ws.Cells.Item(23, 10).Value = 3.4  # J23: 3.45 -> 3.4
$ws.Cells.Item(23, 14).Value = 3.3  # N23: 3.35 -> 3.3
$ws.Cells.Item(23, 16).Value = 1.8  # P23: 1.79 -> 1.8
$ws.Cells.Item(23, 23).Value = 1.73  # W23: 1.74 -> 1.73
# Row 24
$ws.Cells.Item(24, 7).Value = 3  # G24: 3.1 -> 3
$ws.Cells.Item(24, 8).Value = 2.62  # H24: 2.58 -> 2.62
$ws.Cells.Item(24, 9).Value = 2.74  # I24: 2.76 -> 2.74
$ws.Cells.Item(24, 10).Value = 3.5  # J24: 3.45 -> 3.5
$ws.Cells.Item(24, 19).Value = 3.4  # S24: 3.35 -> 3.4
$ws.Cells.Item(24, 22).Value = 1.57  # V24: 1.56 -> 1.57
$ws.Cells.Item(24, 23).Value = 1.5  # W24: 1.48 -> 1.5
# Row 25
$ws.Cells.Item(25, 6).Value = 2.26  # F25: 2.22 -> 2.26
